$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(8)
$sh.TextFrame.TextRange.Text = "08 June 2020"
